$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1, Q1 (copy formatting from O1, then set values)
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update existing data columns for rows 2-25
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# Add new data columns P and Q for rows 2-25
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
